$d = $word.ActiveDocument

# The target sentence lives in its own paragraph:
#   "2. The customer provides a payment method and provides transaction information"
# and must become:
#   "2. The customer enters card info and confirms to pay order"
$p = $d.Paragraphs.Item(13)

# Replace first "provides" -> "enters card info and"
$p.Range.Find.Execute("provides", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "enters card info and", 1) | Out-Null

# Remove " a payment method and "
$p.Range.Find.Execute(" a payment method and ", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "", 1) | Out-Null

# Replace second "provides" -> "confirms"
$p.Range.Find.Execute("provides", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "confirms", 1) | Out-Null

# Replace " transaction information" -> " to pay order"
$p.Range.Find.Execute(" transaction information", $true, $false, $false, $false, $false, `
                 $true, 1, $false, " to pay order", 1) | Out-Null
